$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 1053408  # H19: was 1197032.5
$ws.Cells.Item(19, 9).Value = 1548507.1  # I19: was 1880295.4
$ws.Cells.Item(19, 11).Value = 1548507.1  # K19: was 1880295.4
$ws.Cells.Item(19, 13).Value = -1548332.1  # M19: was -1880120.4
$ws.Cells.Item(111, 8).Value = 835.5  # H111: was 655
$ws.Cells.Item(111, 9).Value = 470  # I111: was 655
$ws.Cells.Item(111, 10).Value = 1932  # J111: was 0
$ws.Cells.Item(111, 11).Value = 1410  # K111: was 1965
$ws.Cells.Item(111, 12).Value = 5796  # L111: was 0
$ws.Cells.Item(111, 13).Value = 1657  # M111: was 1102
$ws.Cells.Item(111, 14).Value = -11930  # N111: was None
$ws.Cells.Item(137, 8).Value = 613569  # H137: was 621535.4
$ws.Cells.Item(137, 9).Value = 1834951.8  # I137: was 1987816.2
$ws.Cells.Item(137, 10).Value = 2877.5962  # J137: was 2842.17
$ws.Cells.Item(137, 11).Value = 5504855.4  # K137: was 5963448.6
$ws.Cells.Item(137, 12).Value = 8632.7886  # L137: was 8526.51
$ws.Cells.Item(137, 13).Value = -5502305.4  # M137: was -5960898.6
$ws.Cells.Item(137, 14).Value = -13732.7886  # N137: was -13626.51
$ws.Cells.Item(138, 8).Value = 3431.6296  # H138: was 3498.291
$ws.Cells.Item(138, 9).Value = 2073.0715  # I138: was 2171
$ws.Cells.Item(138, 10).Value = 3907.125  # J138: was 3909.1191
$ws.Cells.Item(138, 11).Value = 6219.2145  # K138: was 6513
$ws.Cells.Item(138, 12).Value = 11721.375  # L138: was 11727.3573
$ws.Cells.Item(138, 13).Value = -1079.2145  # M138: was -1373
$ws.Cells.Item(138, 14).Value = -22001.375  # N138: was -22007.3573

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(34, 8).Value = 40000  # H34: was 0
$ws.Cells.Item(34, 10).Value = 40000  # J34: was 0
$ws.Cells.Item(34, 12).Value = 40000  # L34: was 0
$ws.Cells.Item(34, 14).Value = -40542  # N34: was None
$ws.Cells.Item(45, 8).Value = 1790.0769  # H45: was 2011.1818
$ws.Cells.Item(45, 9).Value = 2004.6666  # I45: was 2224.125
$ws.Cells.Item(45, 10).Value = 1307.25  # J45: was 1443.3334
$ws.Cells.Item(45, 11).Value = 2004.6666  # K45: was 2224.125
$ws.Cells.Item(45, 12).Value = 1307.25  # L45: was 1443.3334
$ws.Cells.Item(45, 13).Value = -1627.6666  # M45: was -1847.125
$ws.Cells.Item(45, 14).Value = -2061.25  # N45: was -2197.3334
$ws.Cells.Item(74, 8).Value = 4083.5386  # H74: was 4315.5
$ws.Cells.Item(74, 9).Value = 3928.8572  # I74: was 4367
$ws.Cells.Item(74, 11).Value = 3928.8572  # K74: was 4367
$ws.Cells.Item(74, 13).Value = -3054.8572  # M74: was -3493
$ws.Cells.Item(77, 8).Value = 4083.5386  # H77: was 4315.5
$ws.Cells.Item(77, 9).Value = 3928.8572  # I77: was 4367
$ws.Cells.Item(77, 11).Value = 19644.286  # K77: was 21835
$ws.Cells.Item(77, 13).Value = -15276.286  # M77: was -17467
$ws.Cells.Item(109, 8).Value = 37000  # H109: was 32559.092
$ws.Cells.Item(109, 10).Value = 37000  # J109: was 32559.092
$ws.Cells.Item(109, 12).Value = 37000  # L109: was 32559.092
$ws.Cells.Item(109, 14).Value = -39774  # N109: was -35333.092
$ws.Cells.Item(111, 8).Value = 0  # H111: was 60644
$ws.Cells.Item(111, 10).Value = 0  # J111: was 60644
$ws.Cells.Item(111, 12).ClearContents()  # L111: was 60644
$ws.Cells.Item(111, 14).Value = 0  # N111: was -68824
$ws.Cells.Item(132, 8).Value = 1618.7446  # H132: was 1637.9318
$ws.Cells.Item(132, 9).Value = 910.3514  # I132: was 939.17145
$ws.Cells.Item(132, 10).Value = 4239.8  # J132: was 4355.3335
$ws.Cells.Item(132, 11).Value = 2731.0542  # K132: was 2817.51435
$ws.Cells.Item(132, 12).Value = 12719.4  # L132: was 13066.0005
$ws.Cells.Item(132, 13).Value = -201.0542  # M132: was -287.5143500000004
$ws.Cells.Item(132, 14).Value = -17779.4  # N132: was -18126.0005
$ws.Cells.Item(137, 8).Value = 46734  # H137: was 45230
$ws.Cells.Item(137, 10).Value = 46734  # J137: was 45230
$ws.Cells.Item(137, 12).Value = 46734  # L137: was 45230
$ws.Cells.Item(137, 14).Value = -56934  # N137: was -55430

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(59, 8).Value = 49999  # H59: was 0
$ws.Cells.Item(59, 10).Value = 49999  # J59: was 0
$ws.Cells.Item(59, 12).Value = 49999  # L59: was 0
$ws.Cells.Item(59, 14).Value = -51693  # N59: was None
$ws.Cells.Item(134, 8).Value = 3856.6216  # H134: was 4077.2285
$ws.Cells.Item(134, 9).Value = 1388.238  # I134: was 1418.6111
$ws.Cells.Item(134, 10).Value = 7096.375  # J134: was 6892.2354
$ws.Cells.Item(134, 11).Value = 4164.714  # K134: was 4255.8333
$ws.Cells.Item(134, 12).Value = 21289.125  # L134: was 20676.7062
$ws.Cells.Item(134, 13).Value = -1629.714  # M134: was -1720.8333
$ws.Cells.Item(134, 14).Value = -26359.125  # N134: was -25746.7062

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2420.386  # H31: was 2506.6726
$ws.Cells.Item(31, 9).Value = 1086.35  # I31: was 1159.579
$ws.Cells.Item(31, 10).Value = 3141.4866  # J31: was 3217.639
$ws.Cells.Item(31, 11).Value = 1086.35  # K31: was 1159.579
$ws.Cells.Item(31, 12).Value = 3141.4866  # L31: was 3217.639
$ws.Cells.Item(31, 13).Value = -791.3499999999999  # M31: was -864.579
$ws.Cells.Item(31, 14).Value = -3731.4866  # N31: was -3807.639
$ws.Cells.Item(34, 8).Value = 2420.386  # H34: was 2506.6726
$ws.Cells.Item(34, 9).Value = 1086.35  # I34: was 1159.579
$ws.Cells.Item(34, 10).Value = 3141.4866  # J34: was 3217.639
$ws.Cells.Item(34, 11).Value = 1086.35  # K34: was 1159.579
$ws.Cells.Item(34, 12).Value = 3141.4866  # L34: was 3217.639
$ws.Cells.Item(34, 13).Value = -884.3499999999999  # M34: was -957.579
$ws.Cells.Item(34, 14).Value = -3545.4866  # N34: was -3621.639
$ws.Cells.Item(58, 8).Value = 2748.5938  # H58: was 2820.742
$ws.Cells.Item(58, 9).Value = 1534.1072  # I58: was 1571.963
$ws.Cells.Item(58, 11).Value = 1534.1072  # K58: was 1571.963
$ws.Cells.Item(58, 13).Value = -1331.1072  # M58: was -1368.963
$ws.Cells.Item(87, 8).Value = 24398.572  # H87: was 23585.715
$ws.Cells.Item(87, 10).Value = 24398.572  # J87: was 23585.715
$ws.Cells.Item(87, 12).Value = 24398.572  # L87: was 23585.715
$ws.Cells.Item(87, 14).Value = -26770.572  # N87: was -25957.715
$ws.Cells.Item(90, 8).Value = 24398.572  # H90: was 23585.715
$ws.Cells.Item(90, 10).Value = 24398.572  # J90: was 23585.715
$ws.Cells.Item(90, 12).Value = 73195.716  # L90: was 70757.145
$ws.Cells.Item(90, 14).Value = -85051.716  # N90: was -82613.145
$ws.Cells.Item(122, 8).Value = 2049.842  # H122: was 2625
$ws.Cells.Item(122, 9).Value = 957.63635  # I122: was 1320.2
$ws.Cells.Item(122, 10).Value = 3551.625  # J122: was 3440.5
$ws.Cells.Item(122, 11).Value = 2872.90905  # K122: was 3960.6
$ws.Cells.Item(122, 12).Value = 10654.875  # L122: was 10321.5
$ws.Cells.Item(122, 13).Value = -422.9090500000002  # M122: was -1510.6
$ws.Cells.Item(122, 14).Value = -15554.875  # N122: was -15221.5
$ws.Cells.Item(136, 8).Value = 2748.5938  # H136: was 2820.742
$ws.Cells.Item(136, 9).Value = 1534.1072  # I136: was 1571.963
$ws.Cells.Item(136, 11).Value = 4602.321599999999  # K136: was 4715.889
$ws.Cells.Item(136, 13).Value = -2052.321599999999  # M136: was -2165.889
$ws.Cells.Item(139, 8).Value = 39445  # H139: was 39990
$ws.Cells.Item(139, 10).Value = 39445  # J139: was 39990
$ws.Cells.Item(139, 12).Value = 39445  # L139: was 39990
$ws.Cells.Item(139, 14).Value = -49725  # N139: was -50270
$ws.Cells.Item(141, 8).Value = 29580  # H141: was 30742.8
$ws.Cells.Item(141, 10).Value = 29580  # J141: was 30742.8
$ws.Cells.Item(141, 12).Value = 29580  # L141: was 30742.8
$ws.Cells.Item(141, 14).Value = -39940  # N141: was -41102.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(20, 8).Value = 7999.6665  # H20: was 8499.75
$ws.Cells.Item(20, 10).Value = 7999.6665  # J20: was 8499.75
$ws.Cells.Item(20, 12).Value = 23998.9995  # L20: was 25499.25
$ws.Cells.Item(20, 14).Value = -24452.9995  # N20: was -25953.25
$ws.Cells.Item(21, 8).Value = 1983.1111  # H21: was 1550
$ws.Cells.Item(21, 9).Value = 356  # I21: was 1187.5
$ws.Cells.Item(21, 10).Value = 15000  # J21: was 3000
$ws.Cells.Item(21, 11).Value = 1068  # K21: was 3562.5
$ws.Cells.Item(21, 12).Value = 45000  # L21: was 9000
$ws.Cells.Item(21, 13).Value = -895  # M21: was -3389.5
$ws.Cells.Item(21, 14).Value = -45346  # N21: was -9346
$ws.Cells.Item(39, 8).Value = 22301.25  # H39: was 14074.4375
$ws.Cells.Item(39, 9).Value = 0  # I39: was 4999
$ws.Cells.Item(39, 10).Value = 22301.25  # J39: was 14679.467
$ws.Cells.Item(39, 11).Value = 0  # K39: was 14997
$ws.Cells.Item(39, 12).ClearContents()  # L39: was 44038.401
$ws.Cells.Item(39, 13).Value = 66903.75  # M39: was -14703
$ws.Cells.Item(39, 14).Value = -67491.75  # N39: was -44626.401
$ws.Cells.Item(68, 8).Value = 1421.5072  # H68: was 1411.1014
$ws.Cells.Item(68, 9).Value = 1099.7142  # I68: was 1076.1818
$ws.Cells.Item(68, 10).Value = 1562.2916  # J68: was 1567.8723
$ws.Cells.Item(68, 11).Value = 3299.1426  # K68: was 3228.5454
$ws.Cells.Item(68, 12).Value = 4686.8748  # L68: was 4703.6169
$ws.Cells.Item(68, 13).Value = -2488.1426  # M68: was -2417.5454
$ws.Cells.Item(68, 14).Value = -6308.8748  # N68: was -6325.6169
$ws.Cells.Item(71, 8).Value = 1421.5072  # H71: was 1411.1014
$ws.Cells.Item(71, 9).Value = 1099.7142  # I71: was 1076.1818
$ws.Cells.Item(71, 10).Value = 1562.2916  # J71: was 1567.8723
$ws.Cells.Item(71, 11).Value = 9897.427799999999  # K71: was 9685.636200000001
$ws.Cells.Item(71, 12).Value = 14060.6244  # L71: was 14110.8507
$ws.Cells.Item(71, 13).Value = -5841.427799999999  # M71: was -5629.636200000001
$ws.Cells.Item(71, 14).Value = -22172.6244  # N71: was -22222.8507
$ws.Cells.Item(113, 8).Value = 4808311.5  # H113: was 5682444
$ws.Cells.Item(113, 9).Value = 642.06665  # I113: was 657.3333
$ws.Cells.Item(113, 10).Value = 11364224  # J113: was 12500588
$ws.Cells.Item(113, 11).Value = 1926.19995  # K113: was 1971.9999
$ws.Cells.Item(113, 12).Value = 34092672  # L113: was 37501764
$ws.Cells.Item(113, 13).Value = 243.8000500000001  # M113: was 198.0001
$ws.Cells.Item(113, 14).Value = -34097012  # N113: was -37506104
$ws.Cells.Item(131, 8).Value = 907.29  # H131: was 908.1313
$ws.Cells.Item(131, 10).Value = 983.80896  # J131: was 985.625
$ws.Cells.Item(131, 12).Value = 2951.42688  # L131: was 2956.875
$ws.Cells.Item(131, 14).Value = -13031.42688  # N131: was -13036.875
$ws.Cells.Item(132, 8).Value = 2253.6924  # H132: was 2012.85
$ws.Cells.Item(132, 9).Value = 928.5714  # I132: was 960
$ws.Cells.Item(132, 10).Value = 3799.6667  # J132: was 2363.8
$ws.Cells.Item(132, 11).Value = 8357.142600000001  # K132: was 8640
$ws.Cells.Item(132, 12).Value = 34197.0003  # L132: was 21274.2
$ws.Cells.Item(132, 13).Value = -5827.142600000001  # M132: was -6110
$ws.Cells.Item(132, 14).Value = -39257.0003  # N132: was -26334.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 17860422  # H80: was 14709040
$ws.Cells.Item(80, 9).Value = 27781166  # I80: was 20836524
$ws.Cells.Item(80, 11).Value = 27781166  # K80: was 20836524
$ws.Cells.Item(80, 13).Value = -27780168  # M80: was -20835526
$ws.Cells.Item(83, 8).Value = 17860422  # H83: was 14709040
$ws.Cells.Item(83, 9).Value = 27781166  # I83: was 20836524
$ws.Cells.Item(83, 11).Value = 138905830  # K83: was 104182620
$ws.Cells.Item(83, 13).Value = -138900838  # M83: was -104177628
$ws.Cells.Item(113, 8).Value = 2314.3  # H113: was 2157.4546
$ws.Cells.Item(113, 9).Value = 2314.3  # I113: was 2157.4546
$ws.Cells.Item(113, 11).Value = 2314.3  # K113: was 2157.4546
$ws.Cells.Item(113, 13).Value = -144.3000000000002  # M113: was 12.54539999999997
$ws.Cells.Item(132, 8).Value = 4363.375  # H132: was 4972.0557
$ws.Cells.Item(132, 9).Value = 2546.2856  # I132: was 2800
$ws.Cells.Item(132, 10).Value = 5111.5884  # J132: was 5243.5625
$ws.Cells.Item(132, 11).Value = 7638.8568  # K132: was 8400
$ws.Cells.Item(132, 12).Value = 15334.7652  # L132: was 15730.6875
$ws.Cells.Item(132, 13).Value = -5108.8568  # M132: was -5870
$ws.Cells.Item(132, 14).Value = -20394.7652  # N132: was -20790.6875
$ws.Cells.Item(141, 8).Value = 33193  # H141: was 38189.855
$ws.Cells.Item(141, 10).Value = 32937.816  # J141: was 38554.832
$ws.Cells.Item(141, 12).Value = 32937.816  # L141: was 38554.832
$ws.Cells.Item(141, 14).Value = -43297.816  # N141: was -48914.832

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 2529.6667  # H100: was 2410.6365
$ws.Cells.Item(100, 9).Value = 2066.6667  # I100: was 2020
$ws.Cells.Item(100, 10).Value = 2761.1667  # J100: was 2633.8572
$ws.Cells.Item(100, 11).Value = 2066.6667  # K100: was 2020
$ws.Cells.Item(100, 12).Value = 2761.1667  # L100: was 2633.8572
$ws.Cells.Item(100, 13).Value = -1525.6667  # M100: was -1479
$ws.Cells.Item(100, 14).Value = -3843.1667  # N100: was -3715.8572
$ws.Cells.Item(122, 8).Value = 4835  # H122: was 6219.4
$ws.Cells.Item(122, 9).Value = 4326.1113  # I122: was 5209.5
$ws.Cells.Item(122, 10).Value = 7125  # J122: was 10259
$ws.Cells.Item(122, 11).Value = 12978.3339  # K122: was 15628.5
$ws.Cells.Item(122, 12).Value = 21375  # L122: was 30777
$ws.Cells.Item(122, 13).Value = -10528.3339  # M122: was -13178.5
$ws.Cells.Item(122, 14).Value = -26275  # N122: was -35677
$ws.Cells.Item(139, 8).Value = 48683.57  # H139: was 48707.5
$ws.Cells.Item(139, 10).Value = 48683.57  # J139: was 48707.5
$ws.Cells.Item(139, 12).Value = 48683.57  # L139: was 48707.5
$ws.Cells.Item(139, 14).Value = -58963.57  # N139: was -58987.5
$ws.Cells.Item(140, 8).Value = 67809.52  # H140: was 68747.82000000001
$ws.Cells.Item(140, 10).Value = 67809.52  # J140: was 68747.82000000001
$ws.Cells.Item(140, 12).Value = 67809.52  # L140: was 68747.82000000001
$ws.Cells.Item(140, 14).Value = -78169.52  # N140: was -79107.82000000001
$ws.Cells.Item(141, 8).Value = 32216.176  # H141: was 32565.715
$ws.Cells.Item(141, 10).Value = 32216.176  # J141: was 32565.715
$ws.Cells.Item(141, 12).Value = 32216.176  # L141: was 32565.715
$ws.Cells.Item(141, 14).Value = -42576.176  # N141: was -42925.715

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(105, 8).Value = 30123  # H105: was 38750
$ws.Cells.Item(105, 10).Value = 30123  # J105: was 38750
$ws.Cells.Item(105, 12).Value = 30123  # L105: was 38750
$ws.Cells.Item(105, 14).Value = -37111  # N105: was -45738
$ws.Cells.Item(138, 8).Value = 44739.5  # H138: was 41999.332
$ws.Cells.Item(138, 10).Value = 44739.5  # J138: was 41999.332
$ws.Cells.Item(138, 12).Value = 44739.5  # L138: was 41999.332
$ws.Cells.Item(138, 14).Value = -55019.5  # N138: was -52279.332
$ws.Cells.Item(139, 8).Value = 47380  # H139: was 47508
$ws.Cells.Item(139, 10).Value = 47380  # J139: was 47508
$ws.Cells.Item(139, 12).Value = 47380  # L139: was 47508
$ws.Cells.Item(139, 14).Value = -57660  # N139: was -57788
$ws.Cells.Item(140, 8).Value = 42573.375  # H140: was 42948.285
$ws.Cells.Item(140, 10).Value = 51764.5  # J140: was 54127.6
$ws.Cells.Item(140, 12).Value = 51764.5  # L140: was 54127.6
$ws.Cells.Item(140, 14).Value = -62124.5  # N140: was -64487.6
$ws.Cells.Item(141, 8).Value = 48721.668  # H141: was 48723
$ws.Cells.Item(141, 10).Value = 48721.668  # J141: was 48723
$ws.Cells.Item(141, 12).Value = 48721.668  # L141: was 48723
$ws.Cells.Item(141, 14).Value = -59081.668  # N141: was -59083
